# Applies crypto price/volume/coin updates per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.553.20"
$ws.Range("D3").Value = "'1.656.02"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'213.91"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").Value = "'0.509"
$ws.Range("E6").Value = "  -2.43%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'23.93"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -2.41%  "
$ws.Range("D10").Value = "'0.0618"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").Value = "'0.0879"
$ws.Range("D12").Value = "'1.888.58"
$ws.Range("E12").Value = "  -4.30%  "
$ws.Range("D13").Value = "'1.656.85"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").Value = "'0.562"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "'65.79"
$ws.Range("D17").Value = "'27.548.11"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "'240.56"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "'0.0₃0728"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'4.46"
$ws.Range("E22").Value = "  -3.97%  "
$ws.Range("D23").Value = "'9.39"
$ws.Range("E23").Value = "  -3.01%  "
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("D25").Value = "'145.66"
$ws.Range("E25").Value = "  -2.60%  "
$ws.Range("D26").Value = "'7.18"
$ws.Range("E26").Value = "  -4.04%  "
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.21"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0501"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D33").Value = "'1.447.48"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("E34").Value = "  -5.11%  "
$ws.Range("E35").Value = "  -5.22%  "
$ws.Range("E36").Value = "  -1.07%  "
$ws.Range("E37").Value = "  -5.49%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.570"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0171"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("D41").Value = "'69.04"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.22"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.41"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.794"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("D46").Value = "'1.798.69"
$ws.Range("E46").Value = "  -4.29%  "
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "'88.38"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  -6.41%  "
$ws.Range("D51").Value = "'7.82"
$ws.Range("E51").Value = "  -4.22%  "
